# Upgrade DOI prefix 10.5072 -> 10.5281 (Zenodo sandbox -> production prefix)
# across the five occurrences in the deck, preserving all other text/formatting.
# Each target paragraph holds the DOI text in a single run, so we update
# that run's .Text directly (setting .Text on a whole Paragraphs() range
# would otherwise fragment it into multiple runs around the changed chars).

$p = $ppt.ActivePresentation

# Slide 5: "Content Placeholder 2" (shape 2), paragraph 1, run 1.
$s5 = $p.Slides.Item(5)
$run5 = $s5.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1)
$run5.Text = "Museum für Naturkunde. (2024). Photo of Specimen BMT0009388. Zenodo. https://doi.org/10.5281/zenodo.13342373"

# Slide 6: "Content Placeholder 2" (shape 2), paragraph 2, run 1 - the curl command.
$s6 = $p.Slides.Item(6)
$run6 = $s6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2, 1).Runs(1, 1)
$run6.Text = "curl -L ""https://doi.org/10.5281/zenodo.13342373""\
 > bug.tiff "

# Slide 8: "Content Placeholder 2" (shape 2), paragraph 4, run 1.
$s8 = $p.Slides.Item(8)
$run8 = $s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4, 1).Runs(1, 1)
$run8.Text = "load https://doi.org/10.5281/zenodo.13342373"

# Slide 10: "Content Placeholder 2" (shape 2), paragraph 1, run 1.
$s10 = $p.Slides.Item(10)
$run10 = $s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1)
$run10.Text = "Museum für Naturkunde. (2024). Photo of Specimen BMT0009388. Zenodo. https://doi.org/10.5281/zenodo.13342373"

# Slide 11: "Content Placeholder 2" (shape 2), paragraph 3, run 1 (keep run 2 - the bold hash - untouched).
$s11 = $p.Slides.Item(11)
$run11 = $s11.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3, 1).Runs(1, 1)
$run11.Text = "Museum für Naturkunde. (2024). Photo of Specimen BMT0009388. Zenodo. https://doi.org/10.5281/zenodo.13342373 "
